# Student side Activities submission:
# The "CourseDesignerName" (K2) and "LearningCourseName" (L2) values on the
# STAGE sheet are refreshed with newly generated automation values, as part
# of re-running the learning-course/course-designer data submission.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

$ws.Range("K2").Value = "CourseDesigner40081"
$ws.Range("L2").Value = "LearningCourse39214"
